# Add validation to Student in student_import
# The "jezyk" (language) column F is no longer part of the import sheet,
# so remove it entirely and let the remaining columns (grupa, nr tel,
# email, notatka rekrutacyjna) shift one column to the left.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stash the formatting (blue "mailto" font) of the two hyperlink cells in a
# scratch area (far below the real data, in column A so it is unaffected by
# the column shift below) so it can be restored after the hyperlinks are
# recreated further down (Hyperlinks.Add() always repaints the destination
# cell with the built-in "Hyperlink" style, which we don't want here).
$ws.Range("I2").Copy($ws.Range("A100"))
$ws.Range("I3").Copy($ws.Range("A101"))

# Remove column F ("jezyk") completely; G:J shift left to F:I.
$ws.Range("F1:F3").EntireColumn.Delete()

# The hyperlink anchors are not automatically re-targeted by the column
# delete above, so drop the stale ones (still pointing at column I) and
# recreate them at their new home, column H.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("H2"), "mailto:olo@gmail.com", "", "", "olo@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H3"), "mailto:ala@gmail.com", "", "", "ala@gmail.com") | Out-Null

# Restore the original cell formatting that Hyperlinks.Add() just clobbered.
$ws.Range("A100").Copy()
$ws.Range("H2").PasteSpecial(-4122)
$ws.Range("A101").Copy()
$ws.Range("H3").PasteSpecial(-4122)

# Clean up the scratch cells used to stash formatting.
$ws.Range("A100:A101").Clear()

# Match the reviewer's cursor position recorded in the saved file.
$ws.Range("F1").Select()
